$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / tab title to reflect the new "through" date
$ws.Name = "Through 2022-08-26"

# Update the August row label to reflect the new "through" date
$ws.Range("A9").Value = "August (through 08-26)"

# Update August row (row 9) values with the new day's data folded in
$ws.Range("B9").Value = 27
$ws.Range("C9").Value = 63
$ws.Range("D9").Value = 75
$ws.Range("E9").Value = 52
$ws.Range("F9").Value = 38
$ws.Range("G9").Value = 145
$ws.Range("H9").Value = 137
$ws.Range("I9").Value = 139

# Update Total row (row 10) values accordingly
$ws.Range("B10").Value = 189
$ws.Range("C10").Value = 365
$ws.Range("D10").Value = 540
$ws.Range("E10").Value = 477
$ws.Range("F10").Value = 342
$ws.Range("G10").Value = 766
$ws.Range("H10").Value = 1047
$ws.Range("I10").Value = 1110
